# Adds two new columns, I ("I0") and J ("IF"), to the data table.
# Header row (row 1) gets the same bold/bordered/centered style used by
# the other header cells (style index 1, as already applied to H1).
# Data rows 2-60 get the numeric values from the new columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header style (bold, border, centered) from H1 onto
# the two new header cells so they match the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-60) ---
$iVals = @(9,9,8,8,9,8,8,6,4,5,7,6,7,8,6,7,5,8,6,6,6,7,7,7,8,9,8,8,6,7,7,8,5,3,6,7,10,4,6,6,7,6,6,4,6,6,6,8,9,8,7,2,6,6,4,4,6,5,5)
$jVals = @(9,9,9,9,9,8,8,6,5,5,8,6,7,8,6,8,5,9,7,6,6,7,7,7,8,9,8,8,6,7,7,9,5,4,6,7,10,5,7,7,7,6,6,5,7,6,7,8,9,8,7,2,6,6,5,4,6,5,5)

$startRow = 2
for ($k = 0; $k -lt $iVals.Count; $k++) {
    $row = $startRow + $k
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
